$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").Value = "ClosedOn"
$ws.Range("G1").Value = "SubStatus"
$ws.Range("E1").Value = "OrderStatus"
$ws.Range("K1").Value = "InsertServiceRequest"
$ws.Range("K2").Value = "SVMXC__Service_Request__c SR_1 = new SVMXC__Service_Request__c(SVMXC__Status__c = 'Open' );insert SR_1 ;"
$ws.Range("A3").Value = "001q000000hmfgo"
$ws.Range("B3").Value = "a1Jq0000001faAI"
$ws.Columns.Item(11).ColumnWidth = 41
$ws.PageSetup.Orientation = 1
$ws.Range("C13").Select()
